$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read current values (column A holds 8 rows of alternating name/phone entries)
$vals = @()
for ($r = 1; $r -le 8; $r++) {
    $vals += $ws.Cells.Item($r, 1).Value2
}

# Clear the old data region entirely before laying out the new 2-column shape
$ws.Range("A1:A8").ClearContents()

# Rebuild as two columns: column A = phone numbers, column B = names
# Original rows: 1=Name,2=Phone,3=Name,4=Phone,5=Name,6=Phone,7=Name,8=Phone
$destRow = 1
for ($i = 0; $i -lt 8; $i += 2) {
    $name = $vals[$i]
    $phone = $vals[$i + 1]
    $ws.Cells.Item($destRow, 1).Value2 = $phone
    $ws.Cells.Item($destRow, 2).Value2 = $name
    $destRow++
}

# Column B width to match bestFit sizing
$ws.Columns.Item(2).ColumnWidth = 39.7109375

# Update the selection to match the target state
$ws.Range("B9").Select()
